# Update scripts with new TPM: recompute the Ccl11-Ccr5 LR-pair table values
# (rows 2-7) with the new TPM-derived numbers, and drop the rows for the
# self-to-self cluster pairs (e.g. ECs -> ECs) that no longer appear in the
# recomputed output, shrinking the used range from A1:T10 to A1:T7.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.259924
$ws.Range("H2").Value = 0.779772
$ws.Range("I2").Value = 0.0006491957374851489
$ws.Range("J2").Value = 0.0006491957374851488
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.016376
$ws.Range("N2").Value = 0.049128
$ws.Range("O2").Value = 0.4917717717717718
$ws.Range("P2").Value = 0.4917717717717718
$ws.Range("Q2").Value = 0.004256515423999999
$ws.Range("R2").Value = 0.038308638816
$ws.Range("S2").Value = 0.0003192561380497537
$ws.Range("T2").Value = 0.0003192561380497537

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.259924
$ws.Range("H3").Value = 0.779772
$ws.Range("I3").Value = 0.0006491957374851489
$ws.Range("J3").Value = 0.0006491957374851488
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.016924
$ws.Range("N3").Value = 0.050772
$ws.Range("O3").Value = 0.5082282282282282
$ws.Range("P3").Value = 0.5082282282282282
$ws.Range("Q3").Value = 0.004398953775999999
$ws.Range("R3").Value = 0.039590583984
$ws.Range("S3").Value = 0.0003299395994353952
$ws.Range("T3").Value = 0.0003299395994353951

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 383.1307676666667
$ws.Range("H4").Value = 1149.392303
$ws.Range("I4").Value = 0.9569214896224009
$ws.Range("J4").Value = 0.9569214896224006
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.016376
$ws.Range("N4").Value = 0.049128
$ws.Range("O4").Value = 0.4917717717717718
$ws.Range("P4").Value = 0.4917717717717718
$ws.Range("Q4").Value = 6.274149451309333
$ws.Range("R4").Value = 56.467345061784
$ws.Range("S4").Value = 0.4705869763980912
$ws.Range("T4").Value = 0.4705869763980911

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 383.1307676666667
$ws.Range("H5").Value = 1149.392303
$ws.Range("I5").Value = 0.9569214896224009
$ws.Range("J5").Value = 0.9569214896224006
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.016924
$ws.Range("N5").Value = 0.050772
$ws.Range("O5").Value = 0.5082282282282282
$ws.Range("P5").Value = 0.5082282282282282
$ws.Range("Q5").Value = 6.484105111990667
$ws.Range("R5").Value = 58.356946007916
$ws.Range("S5").Value = 0.4863345132243097
$ws.Range("T5").Value = 0.4863345132243095

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 16.98778433333333
$ws.Range("H6").Value = 50.963353
$ws.Range("I6").Value = 0.042429314640114
$ws.Range("J6").Value = 0.04242931464011399
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.016376
$ws.Range("N6").Value = 0.049128
$ws.Range("O6").Value = 0.4917717717717718
$ws.Range("P6").Value = 0.4917717717717718
$ws.Range("Q6").Value = 0.2781919562426666
$ws.Range("R6").Value = 2.503727606184
$ws.Range("S6").Value = 0.02086553923563084
$ws.Range("T6").Value = 0.02086553923563083

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 16.98778433333333
$ws.Range("H7").Value = 50.963353
$ws.Range("I7").Value = 0.042429314640114
$ws.Range("J7").Value = 0.04242931464011399
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.016924
$ws.Range("N7").Value = 0.050772
$ws.Range("O7").Value = 0.5082282282282282
$ws.Range("P7").Value = 0.5082282282282282
$ws.Range("Q7").Value = 0.2875012620573333
$ws.Range("R7").Value = 2.587511358516
$ws.Range("S7").Value = 0.02156377540448316
$ws.Range("T7").Value = 0.02156377540448316

# Remove obsolete rows 8-10 (self-pair rows removed in new TPM data)
$ws.Range("A8:T10").Delete()